$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates: target cell address -> new text value.
# Values that look numeric (e.g. "0.9989") must be force-written as text
# so the saved cell keeps the original text-cell encoding instead of being
# auto-converted into a numeric cell by Excel.
$updates = @(
    @{Cell = 'D2'; Value = '26.061.17'; ForceText = $false}
    @{Cell = 'E2'; Value = '  +0.92%  '; ForceText = $false}
    @{Cell = 'D3'; Value = '1.749.26'; ForceText = $false}
    @{Cell = 'E3'; Value = '  +0.53%  '; ForceText = $false}
    @{Cell = 'D4'; Value = '0.9989'; ForceText = $true}
    @{Cell = 'E4'; Value = '  -0.17%  '; ForceText = $false}
    @{Cell = 'D5'; Value = '233.89'; ForceText = $true}
    @{Cell = 'E5'; Value = '  +3.19%  '; ForceText = $false}
    @{Cell = 'D6'; Value = '0.9992'; ForceText = $true}
    @{Cell = 'E6'; Value = '  -0.10%  '; ForceText = $false}
    @{Cell = 'D7'; Value = '0.5292'; ForceText = $true}
    @{Cell = 'E7'; Value = '  +2.58%  '; ForceText = $false}
    @{Cell = 'D8'; Value = '0.2780'; ForceText = $true}
    @{Cell = 'E8'; Value = '  +2.07%  '; ForceText = $false}
    @{Cell = 'D9'; Value = '0.06194'; ForceText = $true}
    @{Cell = 'E9'; Value = '  +1.73%  '; ForceText = $false}
    @{Cell = 'D10'; Value = '1.749.27'; ForceText = $false}
    @{Cell = 'E10'; Value = '  +0.38%  '; ForceText = $false}
    @{Cell = 'D11'; Value = '0.07256'; ForceText = $true}
    @{Cell = 'E11'; Value = '  +3.61%  '; ForceText = $false}
    @{Cell = 'D12'; Value = '15.34'; ForceText = $true}
    @{Cell = 'E12'; Value = '  +1.08%  '; ForceText = $false}
    @{Cell = 'D13'; Value = '0.6440'; ForceText = $true}
    @{Cell = 'E13'; Value = '  +2.08%  '; ForceText = $false}
    @{Cell = 'D14'; Value = '4.617'; ForceText = $true}
    @{Cell = 'E14'; Value = '  +2.69%  '; ForceText = $false}
    @{Cell = 'D15'; Value = '78.47'; ForceText = $true}
    @{Cell = 'E15'; Value = '  +2.96%  '; ForceText = $false}
    @{Cell = 'D16'; Value = '0.9995'; ForceText = $true}
    @{Cell = 'E16'; Value = '  -0.08%  '; ForceText = $false}
    @{Cell = 'D17'; Value = '0.9992'; ForceText = $true}
    @{Cell = 'E17'; Value = '  -0.09%  '; ForceText = $false}
    @{Cell = 'D18'; Value = '25.977.48'; ForceText = $false}
    @{Cell = 'D19'; Value = '11.62'; ForceText = $true}
    @{Cell = 'E19'; Value = '  +1.60%  '; ForceText = $false}
    @{Cell = 'D20'; Value = '0.000006739'; ForceText = $true}
    @{Cell = 'E20'; Value = '  +2.01%  '; ForceText = $false}
    @{Cell = 'D21'; Value = '1.976.30'; ForceText = $false}
    @{Cell = 'E21'; Value = '  +0.94%  '; ForceText = $false}
    @{Cell = 'D22'; Value = '4.322'; ForceText = $true}
    @{Cell = 'E22'; Value = '  +6.61%  '; ForceText = $false}
    @{Cell = 'D23'; Value = '8.819'; ForceText = $true}
    @{Cell = 'E23'; Value = '  +4.74%  '; ForceText = $false}
    @{Cell = 'D24'; Value = '5.224'; ForceText = $true}
    @{Cell = 'E24'; Value = '  +2.66%  '; ForceText = $false}
    @{Cell = 'D25'; Value = '139.22'; ForceText = $true}
    @{Cell = 'E25'; Value = '  +1.83%  '; ForceText = $false}
    @{Cell = 'D26'; Value = '1.513'; ForceText = $true}
    @{Cell = 'E26'; Value = '  +0.48%  '; ForceText = $false}
    @{Cell = 'D27'; Value = '15.34'; ForceText = $true}
    @{Cell = 'E27'; Value = '  +2.49%  '; ForceText = $false}
    @{Cell = 'D28'; Value = '1.813'; ForceText = $true}
    @{Cell = 'E28'; Value = '  -0.17%  '; ForceText = $false}
    @{Cell = 'D29'; Value = '104.71'; ForceText = $true}
    @{Cell = 'E29'; Value = '  +2.07%  '; ForceText = $false}
    @{Cell = 'D30'; Value = '0.08306'; ForceText = $true}
    @{Cell = 'E30'; Value = '  +0.02%  '; ForceText = $false}
    @{Cell = 'D31'; Value = '3.795'; ForceText = $true}
    @{Cell = 'E31'; Value = '  +4.84%  '; ForceText = $false}
    @{Cell = 'D32'; Value = '3.667'; ForceText = $true}
    @{Cell = 'E32'; Value = '  +8.66%  '; ForceText = $false}
    @{Cell = 'D33'; Value = '0.04531'; ForceText = $true}
    @{Cell = 'E33'; Value = '  +2.93%  '; ForceText = $false}
    @{Cell = 'D34'; Value = '2.642'; ForceText = $true}
    @{Cell = 'E34'; Value = '  +1.29%  '; ForceText = $false}
    @{Cell = 'D35'; Value = '1.003'; ForceText = $true}
    @{Cell = 'E35'; Value = '  +3.67%  '; ForceText = $false}
    @{Cell = 'D36'; Value = '0.6315'; ForceText = $true}
    @{Cell = 'E36'; Value = '  +5.90%  '; ForceText = $false}
    @{Cell = 'D37'; Value = '2.710'; ForceText = $true}
    @{Cell = 'E37'; Value = '  +1.31%  '; ForceText = $false}
    @{Cell = 'D38'; Value = '0.01597'; ForceText = $true}
    @{Cell = 'E38'; Value = '  +2.67%  '; ForceText = $false}
    @{Cell = 'D39'; Value = '1.935'; ForceText = $true}
    @{Cell = 'E39'; Value = '  +0.03%  '; ForceText = $false}
    @{Cell = 'D40'; Value = '0.9986'; ForceText = $true}
    @{Cell = 'E40'; Value = '  -0.07%  '; ForceText = $false}
    @{Cell = 'D41'; Value = '98.24'; ForceText = $true}
    @{Cell = 'E41'; Value = '  -3.10%  '; ForceText = $false}
    @{Cell = 'D42'; Value = '0.3911'; ForceText = $true}
    @{Cell = 'E42'; Value = '  +2.85%  '; ForceText = $false}
    @{Cell = 'D43'; Value = '0.7385'; ForceText = $true}
    @{Cell = 'E43'; Value = '  +2.00%  '; ForceText = $false}
    @{Cell = 'D44'; Value = '5.052'; ForceText = $true}
    @{Cell = 'E44'; Value = '  +3.62%  '; ForceText = $false}
    @{Cell = 'E45'; Value = '  +3.75%  '; ForceText = $false}
    @{Cell = 'D46'; Value = '6.344'; ForceText = $true}
    @{Cell = 'E46'; Value = '  +2.57%  '; ForceText = $false}
    @{Cell = 'D47'; Value = '0.05346'; ForceText = $true}
    @{Cell = 'E47'; Value = '  -2.55%  '; ForceText = $false}
    @{Cell = 'D48'; Value = '54.09'; ForceText = $true}
    @{Cell = 'E48'; Value = '  +4.23%  '; ForceText = $false}
    @{Cell = 'D49'; Value = '30.69'; ForceText = $true}
    @{Cell = 'E49'; Value = '  +3.07%  '; ForceText = $false}
    @{Cell = 'D50'; Value = '7.666'; ForceText = $true}
    @{Cell = 'E50'; Value = '  +3.90%  '; ForceText = $false}
    @{Cell = 'D51'; Value = '0.3467'; ForceText = $true}
    @{Cell = 'E51'; Value = '  +2.37%  '; ForceText = $false}
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.ForceText) {
        $rng.NumberFormat = "@"
    }
    $rng.Value = $u.Value
}
